# Updates the frame-time documentation sentence:
#   " The time between each frame is 0.125 seconds, which means there are 8
#   frames within a second. ..."
# becomes (split across several runs, matching the authored edit):
#   " The time between each frame is 0." | "02" | " seconds, which means
#   there are " | "50" | " frames within a second. ..."

$d = $word.ActiveDocument

# Locate the exact sentence (old text) inside the document body.
$oldText = " The time between each frame is 0.125 seconds, which means there are 8 frames within a second. Despite this constraint, students are still expected to create multiple tasks to solve each problem scenario " + [char]0x2013 + " it is not allowed to define one " + [char]0x201C + "super-task" + [char]0x201D + " that handles all functionality. "

$findRng = $d.Content
$found = $findRng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the frame-time sentence to update."
}

# Re-seat the hit as a brand-new Range: InsertXML only splices cleanly into
# the exact addressed range when that Range wasn't itself produced by
# Find.Execute (the find-hit range's XML insert otherwise lands at the end
# of the story instead of in place).
$rng = $d.Range($findRng.Start, $findRng.End)

# Build the replacement as a run-split fragment (OOXML) and splice it into
# the exact range that was found, leaving the rest of the paragraph intact.
$newRunsXml = '<w:r><w:t xml:space="preserve"> The time between each frame is 0.</w:t></w:r>' +
              '<w:r><w:t>02</w:t></w:r>' +
              '<w:r><w:t xml:space="preserve"> seconds, which means there are </w:t></w:r>' +
              '<w:r><w:t>50</w:t></w:r>' +
              '<w:r><w:t xml:space="preserve"> frames within a second. Despite this constraint, students are still expected to create multiple tasks to solve each problem scenario &#8211; it is not allowed to define one &#8220;super-task&#8221; that handles all functionality. </w:t></w:r>'

$payload = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $newRunsXml + '</w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

$null = $rng.InsertXML($payload)
